# Apply updated price/volume data from the Jan 25 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.91%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.03%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.060"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.93%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08008"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.85%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.933"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-10.31%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.050"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.28%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.749"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.43%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9213"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.68%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1220"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "21.14%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1854"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.76%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09427"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.81%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03582"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.35%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09842"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.85%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001387"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.29%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005807"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.20%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.490"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.70%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3408"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.90%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1283"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.52%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.040"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.93%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.37%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04524"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.65%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.50%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004845"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.35%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001251"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-7.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01934"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.56%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04755"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.37%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007535"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.50%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009557"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "21.95%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1331"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.81%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002111"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.31%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01099"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.69%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006292"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.59%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "58.62%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001488"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-21.85%"
